$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.201286128529318; "C"=0.3356557433913565; "E"=0.4255167328346232; "F"=0.4443680307746121; "G"=0.4243298004187892; "H"=0.5027521521928691 }
    3 = @{ "B"=1.053283970701386; "C"=0.2944854172539522; "E"=0.3708373693408618; "F"=0.3878228170618172; "G"=0.4081752012819067; "H"=0.5030086175844133 }
    4 = @{ "B"=0.9624234431025798; "C"=0.2691076642603889; "E"=0.3373969496919926; "F"=0.3531389305169483; "G"=0.3990250212695656; "H"=0.5037832496148127 }
    5 = @{ "B"=0.9253994836014385; "C"=0.2587408941710976; "E"=0.3238001751638819; "F"=0.3390132514313251; "G"=0.39548596248693; "H"=0.5042525617249822 }
    6 = @{ "B"=0.9192518053934009; "C"=0.2570179757479139; "E"=0.3215441955878049; "F"=0.336668177824194; "G"=0.3949096573982871; "H"=0.504339725194356 }
    7 = @{ "B"=0.9619241163084098; "C"=0.2689679564061009; "E"=0.3372134593206795; "F"=0.3529483938344953; "G"=0.398976528911021; "H"=0.5037889588845417 }
    8 = @{ "B"=1.150251581965676; "C"=0.3214806048942762; "E"=0.4066338917371013; "F"=0.4248636149813478; "G"=0.4185979467954439; "H"=0.5027115549025325 }
    9 = @{ "B"=1.519730955940133; "C"=0.4236886130502171; "E"=0.5439627685100561; "F"=0.5661985755041457; "G"=0.463340775190062; "H"=0.505566334961685 }
    10 = @{ "B"=1.791410350719502; "C"=0.4983432279729527; "E"=0.6457930840560522; "F"=0.6702781546542269; "G"=0.500264868217414; "H"=0.5107912738788514 }
    11 = @{ "B"=1.915081600282861; "C"=0.532218142210013; "E"=0.6923671657365134; "F"=0.7176906081379002; "G"=0.5179929272211155; "H"=0.5138689993109722 }
    12 = @{ "B"=1.96192646907997; "C"=0.545033839025109; "E"=0.7100432713122729; "F"=0.7356546913071611; "G"=0.5248440704440895; "H"=0.5151370067448227 }
    13 = @{ "B"=1.951836977070855; "C"=0.5422742772429388; "E"=0.7062345931862382; "F"=0.7317853510981394; "G"=0.5233623560704075; "H"=0.514859330039144 }
    14 = @{ "B"=1.918935282368579; "C"=0.533272735682317; "E"=0.693820577850417; "F"=0.7191683204515869; "G"=0.5185537877710544; "H"=0.5139712538843924 }
    15 = @{ "B"=1.898783810783641; "C"=0.5277574797835314; "E"=0.68622188337892; "F"=0.7114413442032514; "G"=0.515626481026402; "H"=0.5134406876677531 }
    16 = @{ "B"=1.783329923805468; "C"=0.4961277242689448; "E"=0.6427546993909061; "F"=0.6671810134426437; "G"=0.4991253768858144; "H"=0.510604390254656 }
    17 = @{ "B"=1.712524859856956; "C"=0.4767021452310587; "E"=0.6161556817610716; "F"=0.6400460337215605; "G"=0.4892438108791453; "H"=0.5090451672384688 }
    18 = @{ "B"=1.671807380953339; "C"=0.4655209635517963; "E"=0.600880043088651; "F"=0.6244449056556647; "G"=0.4836475927953643; "H"=0.5082142070544933 }
    19 = @{ "B"=1.658022429374455; "C"=0.4617338058338305; "E"=0.5957118922657116; "F"=0.6191636801734006; "G"=0.4817676981824803; "H"=0.5079441171771748 }
    20 = @{ "B"=1.720061371349175; "C"=0.4787708674225541; "E"=0.6189847448771815; "F"=0.6429339538360921; "G"=0.4902866461563917; "H"=0.50920431990167 }
    21 = @{ "B"=1.928598937042068; "C"=0.535917029375014; "E"=0.6974657751179762; "F"=0.7228739723492197; "G"=0.5199624057202925; "H"=0.5142293064062642 }
    22 = @{ "B"=2.064968203013393; "C"=0.5731955951446253; "E"=0.7489896073408175; "F"=0.7751780083420101; "G"=0.5401629800573744; "H"=0.5181119324109886 }
    23 = @{ "B"=1.992177769911962; "C"=0.5533055844457522; "E"=0.7214679965777293; "F"=0.7472568307916134; "G"=0.529306488857884; "H"=0.5159843538336588 }
    24 = @{ "B"=1.716654147943757; "C"=0.4778356391277043; "E"=0.6177056743350704; "F"=0.6416283278902171; "G"=0.489814916417771; "H"=0.5091321632287418 }
    25 = @{ "B"=1.419746781554693; "C"=0.3961178371353071; "E"=0.5066626035542185; "F"=0.5279251897347308; "G"=0.39548596248693; "H"=0.5042520343596237 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}